# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) figures
# as scraped on Mon Apr  1 15:52:59 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'" + "68.712.00"
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").Value = "'" + "3.481.50"
$ws.Range("E3").Value = "  -3.78%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'" + "576.11"
$ws.Range("E5").Value = "  -4.29%  "

$ws.Range("D6").Value = "'" + "189.03"
$ws.Range("E6").Value = "  -3.35%  "

$ws.Range("D7").Value = "'" + "3.469.97"
$ws.Range("E7").Value = "  -3.79%  "

$ws.Range("D8").Value = "'" + "0.603"

$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("E10").Value = "  -5.08%  "

$ws.Range("E11").Value = "  -5.12%  "

$ws.Range("D12").Value = "'" + "51.67"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("D13").Value = "'" + "0.0000284"
$ws.Range("E13").Value = "  -6.72%  "

$ws.Range("D14").Value = "'" + "9.08"
$ws.Range("E14").Value = "  -5.13%  "

$ws.Range("D15").Value = "'" + "4.039.64"
$ws.Range("E15").Value = "  -3.59%  "

$ws.Range("D16").Value = "'" + "633.59"
$ws.Range("E16").Value = "  +5.71%  "

$ws.Range("D17").Value = "'" + "68.652.28"
$ws.Range("E17").Value = "  -2.52%  "

$ws.Range("D18").Value = "'" + "3.486.18"
$ws.Range("E18").Value = "  -3.67%  "

$ws.Range("E19").Value = "  -4.44%  "

$ws.Range("E20").Value = "  -2.34%  "

$ws.Range("E21").Value = "  -5.27%  "

$ws.Range("D22").Value = "'" + "0.942"
$ws.Range("E22").Value = "  -5.69%  "

$ws.Range("E23").Value = "  -4.03%  "

$ws.Range("D24").Value = "'" + "5.38"
$ws.Range("E24").Value = "  +3.28%  "

$ws.Range("D25").Value = "'" + "99.12"
$ws.Range("E25").Value = "  -3.46%  "

$ws.Range("E26").Value = "  -6.91%  "

$ws.Range("E27").Value = "  -4.70%  "

$ws.Range("E28").Value = "  +2.15%  "

$ws.Range("D29").Value = "'" + "10.02"
$ws.Range("E29").Value = "  -5.52%  "

$ws.Range("D30").Value = "'" + "9.19"
$ws.Range("E30").Value = "  -5.33%  "

$ws.Range("D31").Value = "'" + "32.37"
$ws.Range("E31").Value = "  -4.21%  "

$ws.Range("D32").Value = "'" + "6.70"
$ws.Range("E32").Value = "  -8.14%  "

$ws.Range("D33").Value = "'" + "4.05"
$ws.Range("E33").Value = "  -14.34%  "

$ws.Range("D34").Value = "'" + "11.58"
$ws.Range("E34").Value = "  -5.69%  "

$ws.Range("E35").Value = "  -7.52%  "

$ws.Range("D37").Value = "'" + "3.700.51"
$ws.Range("E37").Value = "  -6.03%  "

$ws.Range("D38").Value = "'" + "0.999"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").Value = "'" + "0.0₃0786"
$ws.Range("E39").Value = "  -10.99%  "

$ws.Range("D40").Value = "'" + "498.74"
$ws.Range("E40").Value = "  -5.48%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").Value = "'" + "2.92"
$ws.Range("E42").Value = "  -4.41%  "

$ws.Range("E43").Value = "  -5.61%  "

$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").Value = "'" + "34.24"
$ws.Range("E45").Value = "  -7.32%  "

$ws.Range("E46").Value = "  -4.70%  "

$ws.Range("E47").Value = "  -5.51%  "

$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("E49").Value = "  -4.20%  "

$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("D51").Value = "'" + "8.06"
$ws.Range("E51").Value = "  -5.91%  "
